$d = $word.ActiveDocument

# --- 1. "//Naam nog bedenken" -> "CALGON", and move the _GoBack bookmark
#        from the old "Ruimte" paragraph to the end of this (now first) paragraph ---

# Use a temporary trailing marker character "X" so the insertion point that will
# host the (collapsed) bookmark is a genuine mid-paragraph position rather than
# the paragraph-mark boundary; then strip the marker back out.
$d.Content.Find.Execute("//Naam nog bedenken", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CALGONX", 2)

$r = $d.Content
$r.Find.Execute("CALGONX")
$markerPos = $r.End - 1
$bmRange = $d.Range($markerPos, $markerPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Text = ""

# --- 2. "Ruimte: //nog niet bekend" -> "Ruimte: " + "D08.27" (two runs) ---
$d.Content.Find.Execute("Ruimte: //nog niet bekend", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ruimte: D08.27", 2)

$rRuimte = $d.Content
$rRuimte.Find.Execute("D08.27")
$rRuimte.Bold = 1
$rRuimte.Bold = 0

# --- 3. "14:30" -> "13" + ":30" (two runs, effectively becomes "13:30") ---
$rTijd = $d.Content
$rTijd.Find.Execute("14:30")
$rTijdFirst = $d.Range($rTijd.Start, $rTijd.Start + 2)
$rTijdFirst.Text = "13"
$rTijdFirst.Bold = 1
$rTijdFirst.Bold = 0

# --- 4. Insert new list item "Gedane zaken" after
#        "Goedkeuring van notulen van vorige vergadering" ---
$pGoedkeuring = $d.Paragraphs.Item(10)
$pGoedkeuring.Range.InsertParagraphAfter()
$pGedane = $d.Paragraphs.Item(11)
$pGedane.Range.Text = "Gedane zaken"

Write-Output "done"
